# Add newly-added health facilities to the "facilities" choice list (choices sheet)
# and renumber the existing facility labels with their list position, matching the
# commit "add new facilities to xlsx files".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

# --- 1. Prefix the existing facility labels with their numeric position ---
$ws.Range("C5").Value2  = "1. Chitungwiza-Seke North clinic"
$ws.Range("C6").Value2  = "2. Chitungwiza-Seke South clinic"
$ws.Range("C7").Value2  = "3. City Med hospital"
$ws.Range("C8").Value2  = "4. Zengeza Clinic"
$ws.Range("C9").Value2  = "5. Chitungwiza Central Hospital"
$ws.Range("C10").Value2 = "6. Chegutu- Norton hospital"
$ws.Range("C11").Value2 = "7. Chegutu District Hospital"
$ws.Range("C12").Value2 = "8. Monera clinic(Norton Outreach)"
$ws.Range("C13").Value2 = "9. Marondera District Hospital"
$ws.Range("C14").Value2 = "10. Mahusekwa Hospital"

# --- 2. Append the new facilities (rows 15-23), matching the field entry order
#        used by the original author (name/label typed in a slightly mixed
#        order for a couple of rows) so the shared-string table lines up.
$ws.Range("A15").Value2 = "facilities"
$ws.Range("B15").Value2 = "makumbe"
$ws.Range("C15").Value2 = "11. Goromonzi-Makumbe Mission Hospital"

$ws.Range("A16").Value2 = "facilities"
$ws.Range("B16").Value2 = "ruwa"
$ws.Range("C16").Value2 = "12. Goromonzi-Ruwa Rehab Hospital"

$ws.Range("A17").Value2 = "facilities"
$ws.Range("C17").Value2 = "13. Sanyati-Kadoma Hospital"
$ws.Range("B17").Value2 = "kadoma"

$ws.Range("A18").Value2 = "facilities"
$ws.Range("A19").Value2 = "facilities"
$ws.Range("B19").Value2 = "musiso"
$ws.Range("C18").Value2 = "14. Zaka-Ndanga District Hospital"
$ws.Range("C19").Value2 = "15. Zaka-Musiso Mission Hospital"
$ws.Range("B18").Value2 = "ndanga"

$ws.Range("A20").Value2 = "facilities"
$ws.Range("B20").Value2 = "musiso"
$ws.Range("C20").Value2 = "16. Mberengwa-Musiso Mission Hospital"

$ws.Range("A21").Value2 = "facilities"
$ws.Range("B21").Value2 = "musiso"
$ws.Range("C21").Value2 = "17. Mberengwa-Mnene Mission Hospital"

$ws.Range("A22").Value2 = "facilities"
$ws.Range("B22").Value2 = "musiso"
$ws.Range("C22").Value2 = "18. Mberengwa-Musume Mission Hospital"

$ws.Range("A23").Value2 = "facilities"
$ws.Range("B23").Value2 = "musiso"
$ws.Range("C23").Value2 = "19. Mberengwa-Mberengwa District Hospital"

# --- 3. Give the new rows (15-23) the same look as the existing facility rows ---
$ws.Range("A5").Copy() | Out-Null
$ws.Range("A15:A24").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A5").Copy() | Out-Null
$ws.Range("B15:C23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("B15:C23").WrapText = $false

# --- 4. Row 24 is a leftover, formatted-but-empty row below the new entries ---
$ws.Range("A24").Value2 = "facilities"
$ws.Range("A24").ClearContents() | Out-Null

$excel.CutCopyMode = 0
Write-Host "facilities list updated"
